$wb = $excel.ActiveWorkbook

# Rename the worksheets (drop the trailing " Smell" suffix)
$wsLongMethod = $wb.Worksheets.Item("Long Method Smell")
$wsLongMethod.Name = "Long Method"

$wsLargeClass = $wb.Worksheets.Item("Large Class Smell")
$wsLargeClass.Name = "Large Class"

# Move the active selection on the "Long Method" sheet to E17, matching the
# cursor position left behind by the author's editing session.
$wsLongMethod.Activate()
$wsLongMethod.Range("E17").Select()
